$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Update the SQL Server IP address in E2
$ws.Range("E2").Value = "192.168.0.24"

# Move the active selection from G6 to H6 (matches the saved view state)
$ws.Range("H6").Select()
